# Fix IPMI sensor number text: remove the erroneous "sensor number, " prefix
# that appeared starting at DIMM_C1_CPU1 (row 16) through row 62.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: "sensor number, UNC, UC" -> "UNC, UC"
$ws.Range("B16").Value = "UNC, UC"

# Rows 17-47 and 50-61: remove the stray "sensor number, " cell entirely.
$rowsToClear = @(17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,50,51,52,53,54,55,56,57,58,59,60,61)
foreach ($r in $rowsToClear) {
    $ws.Cells.Item($r, 2).ClearContents()
}

# Rows 48-49: "Sensor name, sensor number, " -> "Sensor name, "
$ws.Range("B48").Value = "Sensor name, "
$ws.Range("B49").Value = "Sensor name, "

# Row 62: remove the stray "sensor number, " text, leaving the cell blank.
$ws.Range("B62").ClearContents()
